$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column string values (which may look numeric, e.g. "5.24" or
# "62.702.34") are written as text, not auto-converted to numbers, while
# keeping the cell style identical to the original (unstyled) cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "62.702.34"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "2.436.61"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "566.38"
$ws.Range("D6").Value = "145.33"
$ws.Range("E6").Value = "  +2.06%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("E9").Value = "  +1.99%  "
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("D11").Value = "5.24"
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "26.80"
$ws.Range("E13").Value = "  +4.91%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "0.0000185"
$ws.Range("E14").Value = "  +6.46%  "
$ws.Range("D15").Value = "2.877.54"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "62.474.39"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "2.434.95"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").Value = "11.24"
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("D20").Value = "323.97"
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D23").Value = "67.21"
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("D24").Value = "1.78"
$ws.Range("E24").Value = "  +3.66%  "
$ws.Range("D25").Value = "8.74"
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("E26").Value = "  +8.83%  "
$ws.Range("D27").Value = "566.96"
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("D28").Value = "2.557.78"
$ws.Range("E28").Value = "  +1.18%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").Value = "8.40"
$ws.Range("E30").Value = "  +2.90%  "
$ws.Range("E31").Value = "  +3.01%  "
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("D34").Value = "1.54"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").Value = "4.87"
$ws.Range("E35").Value = "  +4.11%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("D38").Value = "5.43"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("D39").Value = "18.77"
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.84"
$ws.Range("E40").Value = "  +2.59%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "148.38"
$ws.Range("E41").Value = "  -2.18%  "
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("E43").Value = "  +6.75%  "
$ws.Range("D44").Value = "148.53"
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("D45").Value = "3.67"
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("D46").Value = "0.0538"
$ws.Range("E46").Value = "  +1.19%  "
$ws.Range("D47").Value = "20.54"
$ws.Range("E47").Value = "  +3.55%  "
$ws.Range("D48").Value = "0.601"
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "0.0928"
$ws.Range("E49").Value = "  +1.34%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0231"
$ws.Range("E50").Value = "  +2.67%  "
$ws.Range("E51").Value = "  +0.58%  "

# Restore the default "Normal" style on the D column so the number-format
# override above does not leave a lingering style difference.
$ws.Range("D2:D51").Style = "Normal"
